# Reorder the "Recorded By" (column G) entries on the "Session Analysis
# Results" sheet so that email-like entries (containing "@") come first,
# followed by the remaining entries (e.g. "System", "system"), preserving
# the relative order within each group.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ", "

    $emails = @()
    $others = @()
    foreach ($p in $parts) {
        if ($p -like "*@*") {
            $emails += $p
        } else {
            $others += $p
        }
    }

    $newParts = $emails + $others
    $newVal = [string]::Join(", ", $newParts)

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
